# Update column G ("K" - strikeouts) values on Sheet1 to reflect
# regenerated save_data (commit: "regen save_data to use K instead of
# Strike#, regen std/mean, calc and write s_vals").
#
# Only the numeric values in column G, rows 2-31, change; everything
# else in the sheet stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(6, 8, 5, 6, 1, 4, 6, 7, 7, 6, 3, 10, 12, 5, 5, 9, 5, 4, 5, 5, 4, 8, 3, 9, 2, 1, 3, 1, 4, 2)

$row = 2
foreach ($val in $kValues) {
    $ws.Cells.Item($row, 7).Value = $val
    $row++
}
